$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The existing footer block (rows 60-64: two blank rows, a blank separator
# row, the "TC/Verity1.0" label row and the "user/password" sample row) needs
# to move down to make room for:
#   - 7 new regression-test rows (DEC_0874, DEC_0876, DEC_0877, DEC_0879,
#     DEC_0880, DEC_0884, DEC_0885)
#   - 1 new "Verity3.0 / verity" class marker row (DEC_1139)
# Net effect: old row 60 -> new row 67, old rows 61-64 -> new rows 69-72,
# with the brand-new DEC_1139 row landing at 68.
# ---------------------------------------------------------------------------

# Step 1: make room for the 7 new test-case rows above the old block.
$ws.Rows("60:66").Insert()

# Step 2: make room for the new DEC_1139/Verity3.0 marker row, right after
# the first (now shifted) blank separator row.
$ws.Rows("68").Insert()

# ---------------------------------------------------------------------------
# Fill the 7 new regression-test rows (60-66). Each follows the same layout
# as all the other test rows above them: TC id / fixed user / fixed password
# / SIN_DATO placeholders for the remaining BUSQUEDA..EMAIL_REPRESENTANTE
# columns.
# ---------------------------------------------------------------------------
$testIds = @("DEC_0874", "DEC_0876", "DEC_0877", "DEC_0879", "DEC_0880", "DEC_0884", "DEC_0885")
$row = 60
foreach ($testId in $testIds) {
    $ws.Range("A" + $row).Value = $testId
    $ws.Range("B" + $row).Value = "18092588-0"
    $ws.Range("C" + $row).Value = "sebA`$1357"
    $ws.Range("D" + $row + ":J" + $row).Value = "SIN_DATO"
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# New marker row 68: DEC_1139 / 13712759-8 / Verity3.0 / verity, followed by
# SIN_DATO placeholders out to column O.
# ---------------------------------------------------------------------------
$ws.Range("A68").Value = "DEC_1139"
$ws.Range("B68").Value = "13712759-8"
$ws.Range("C68").Value = "Verity3.0"
$ws.Range("D68").Value = "verity"
$ws.Range("E68:O68").Value = "SIN_DATO"

# Restore the view/selection state recorded for this edit.
$ws.Range("H10").Select()
